$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for both rows
$wsOverview.Range("G2").Value = "2016-10-26 08:32:42"
$wsOverview.Range("G3").Value = "2016-10-26 08:32:42"

# zh-cn sheet: Priority "ht" -> "mt" for both rows
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H) and Correspond Handback DateTime (K)
$wsZhCn.Range("H2").Value = "2016-10-26 08:32:30"
$wsZhCn.Range("H3").Value = "2016-10-26 08:32:30"
$wsZhCn.Range("K2").Value = "2016-10-26 08:33:12"
$wsZhCn.Range("K3").Value = "2016-10-26 08:33:12"

# de-de sheet: Priority "ht" -> "mt" for both rows
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# de-de sheet: Correspond Handoff Datetime (H) - matches Overview's updated timestamp
$wsDeDe.Range("H2").Value = "2016-10-26 08:32:42"
$wsDeDe.Range("H3").Value = "2016-10-26 08:32:42"

# de-de sheet: Correspond Handback DateTime (K)
$wsDeDe.Range("K2").Value = "2016-10-26 08:33:29"
$wsDeDe.Range("K3").Value = "2016-10-26 08:33:29"
